$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F8").Value = 43
$ws1.Range("F10").Value = 2702
$ws1.Range("F12").Value = 1796
$ws1.Range("F13").Value = 617
$ws1.Range("F14").Value = 292
$ws1.Range("F15").Value = 704
$ws1.Range("F16").Value = 5169
$ws1.Range("F18").Value = 85
$ws1.Range("F20").Value = 3391
$ws1.Range("F21").Value = 877
$ws1.Range("F25").Value = 2449
$ws1.Range("F30").Value = 492
$ws1.Range("F31").Value = 1311
$ws1.Range("F33").Value = 10
$ws1.Range("F34").Value = 75
$ws1.Range("F37").Value = 1469
$ws1.Range("F38").Value = 25
$ws1.Range("F39").Value = 1425
$ws1.Range("F40").Value = 93

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 19
$ws2.Range("G10").Value = 388
$ws2.Range("F11").Value = 155
$ws2.Range("F18").Value = 264
$ws2.Range("F19").Value = 520

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 905
$ws3.Range("F6").Value = 43
$ws3.Range("F7").Value = 66

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 905
$ws4.Range("F7").Value = 43
$ws4.Range("F8").Value = 66
$ws4.Range("F13").Value = 19
$ws4.Range("F15").Value = 43
$ws4.Range("F16").Value = 2702
$ws4.Range("F20").Value = 1796
$ws4.Range("F21").Value = 155
$ws4.Range("F22").Value = 617
$ws4.Range("F23").Value = 292
$ws4.Range("F24").Value = 704
$ws4.Range("F25").Value = 5169
$ws4.Range("F27").Value = 85
$ws4.Range("F29").Value = 3391
$ws4.Range("F30").Value = 877
$ws4.Range("F35").Value = 2449
$ws4.Range("F38").Value = 492
$ws4.Range("F39").Value = 1311
$ws4.Range("F41").Value = 264
$ws4.Range("F42").Value = 520
$ws4.Range("F44").Value = 10
$ws4.Range("F45").Value = 75
$ws4.Range("F48").Value = 25
$ws4.Range("F50").Value = 1425
$ws4.Range("F51").Value = 93
